$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column C header is not set (no C1 value per diff), but column width is set.
$ws.Columns.Item(3).ColumnWidth = 41

# Fill in the new column C values for the relevant rows
$ws.Range("C4").Value  = "Comparativo 2008-2012 X 2013-2023 (marco Lei das Cotas)"
$ws.Range("C5").Value  = "Inicialmente PISM e SISU (vestibular para o período de 2008-2012)"
$ws.Range("C6").Value  = 'Tirar os "cursos" ABI'
$ws.Range("C8").Value  = "Evasão por tipo de ingresso"
$ws.Range("C9").Value  = "Evasão por cota"
$ws.Range("C10").Value = "Evasão por sexo"
$ws.Range("C11").Value = "Evasão por curso"

# Update the selection to match the post-edit state
$ws.Range("C16").Select()
